# Update LR-pair (Jag2-Notch1) NATMI metrics with newly recomputed TPM-based values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 15.93194
$ws.Range("H2").Value = 47.79582
$ws.Range("I2").Value = 0.9552847657129105
$ws.Range("J2").Value = 0.9552847657129107
$ws.Range("M2").Value = 38.55267666666666
$ws.Range("N2").Value = 115.65803
$ws.Range("O2").Value = 0.5758151725879548
$ws.Range("P2").Value = 0.5758151725879548
$ws.Range("Q2").Value = 614.2189314927332
$ws.Range("R2").Value = 5527.9703834346
$ws.Range("S2").Value = 0.5500674622396236
$ws.Range("T2").Value = 0.5500674622396236
$ws.Range("G3").Value = 15.93194
$ws.Range("H3").Value = 47.79582
$ws.Range("I3").Value = 0.9552847657129105
$ws.Range("J3").Value = 0.9552847657129107
$ws.Range("O3").Value = 0.08021535714867321
$ws.Range("P3").Value = 0.08021535714867323
$ws.Range("Q3").Value = 85.56528779144
$ws.Range("R3").Value = 770.0875901229601
$ws.Range("S3").Value = 0.07662850866034773
$ws.Range("T3").Value = 0.07662850866034776
$ws.Range("G4").Value = 15.93194
$ws.Range("H4").Value = 47.79582
$ws.Range("I4").Value = 0.9552847657129105
$ws.Range("J4").Value = 0.9552847657129107
$ws.Range("M4").Value = 23.02986166666667
$ws.Range("N4").Value = 69.089585
$ws.Range("O4").Value = 0.3439694702633719
$ws.Range("P4").Value = 0.3439694702633719
$ws.Range("Q4").Value = 366.9103742816333
$ws.Range("R4").Value = 3302.1933685347
$ws.Range("S4").Value = 0.3285887948129392
$ws.Range("T4").Value = 0.3285887948129392
$ws.Range("I5").Value = 0.004609931913019111
$ws.Range("J5").Value = 0.004609931913019112
$ws.Range("M5").Value = 38.55267666666666
$ws.Range("N5").Value = 115.65803
$ws.Range("O5").Value = 0.5758151725879548
$ws.Range("P5").Value = 0.5758151725879548
$ws.Range("Q5").Value = 2.964045440163333
$ws.Range("R5").Value = 26.67640896147
$ws.Range("S5").Value = 0.00265446874011382
$ws.Range("T5").Value = 0.002654468740113821
$ws.Range("I6").Value = 0.004609931913019111
$ws.Range("J6").Value = 0.004609931913019112
$ws.Range("O6").Value = 0.08021535714867321
$ws.Range("P6").Value = 0.08021535714867323
$ws.Range("Q6").Value = 0.412913682908
$ws.Range("R6").Value = 3.716223146172001
$ws.Range("S6").Value = 0.0003697873348338943
$ws.Range("T6").Value = 0.0003697873348338944
$ws.Range("I7").Value = 0.004609931913019111
$ws.Range("J7").Value = 0.004609931913019112
$ws.Range("M7").Value = 23.02986166666667
$ws.Range("N7").Value = 69.089585
$ws.Range("O7").Value = 0.3439694702633719
$ws.Range("P7").Value = 0.3439694702633719
$ws.Range("Q7").Value = 1.770604854518333
$ws.Range("R7").Value = 15.935443690665
$ws.Range("S7").Value = 0.001585675838071396
$ws.Range("T7").Value = 0.001585675838071397
$ws.Range("G8").Value = 0.6688636666666667
$ws.Range("H8").Value = 2.006591
$ws.Range("I8").Value = 0.04010530237407027
$ws.Range("J8").Value = 0.04010530237407027
$ws.Range("M8").Value = 38.55267666666666
$ws.Range("N8").Value = 115.65803
$ws.Range("O8").Value = 0.5758151725879548
$ws.Range("P8").Value = 0.5758151725879548
$ws.Range("Q8").Value = 25.78648467508111
$ws.Range("R8").Value = 232.07836207573
$ws.Range("S8").Value = 0.02309324160821739
$ws.Range("T8").Value = 0.02309324160821739
$ws.Range("G9").Value = 0.6688636666666667
$ws.Range("H9").Value = 2.006591
$ws.Range("I9").Value = 0.04010530237407027
$ws.Range("J9").Value = 0.04010530237407027
$ws.Range("O9").Value = 0.08021535714867321
$ws.Range("P9").Value = 0.08021535714867323
$ws.Range("Q9").Value = 3.592250041838668
$ws.Range("R9").Value = 32.33025037654801
$ws.Range("S9").Value = 0.003217061153491579
$ws.Range("T9").Value = 0.003217061153491579
$ws.Range("G10").Value = 0.6688636666666667
$ws.Range("H10").Value = 2.006591
$ws.Range("I10").Value = 0.04010530237407027
$ws.Range("J10").Value = 0.04010530237407027
$ws.Range("M10").Value = 23.02986166666667
$ws.Range("N10").Value = 69.089585
$ws.Range("O10").Value = 0.3439694702633719
$ws.Range("P10").Value = 0.3439694702633719
$ws.Range("Q10").Value = 15.40383771719278
$ws.Range("R10").Value = 138.634539454735
$ws.Range("S10").Value = 0.0137949996123613
$ws.Range("T10").Value = 0.0137949996123613